$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 corresponds to file a4d143dc-a122-4c5d-aea4-1d1a9c87264c.md
# Column D = "Latest Handoff Datetime"
$wsZhCn.Range("D7").Value = "2016-03-10 06:55:38"
$wsDeDe.Range("D7").Value = "2016-03-10 06:55:44"
